$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Cape Schank" -> "Cape Schanck"
$ws.Range("A10").Value = "Cape Schanck"

# Insert a new row before row 41 (old row 41 "Southern Cross..." shifts to 42,
# old row 42 "Wonthaggi..." shifts to 43), and populate the new row with the
# Ringwood exposure site entry.
$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = "Ringwood"
$ws.Range("B41").Value = "Block 7 Dumplings, 171 - 175 Maroondah Highway"
$ws.Range("C41").Value = "29/12/20 8:10pm - 8:45pm"
$ws.Range("D41").Value = "Case attended store"
